$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values that look numeric stay as literal text (matches source formatting).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '42.652.44'
$ws.Range("E2").Value = '  -1.56%  '
$ws.Range("D3").Value = '2.304.19'
$ws.Range("E3").Value = '  +0.04%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '303.26'
$ws.Range("E5").Value = '  -1.99%  '
$ws.Range("D6").Value = '99.52'
$ws.Range("E6").Value = '  -4.79%  '
$ws.Range("E7").Value = '  -4.34%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").Value = '0.503'
$ws.Range("E9").Value = '  -4.36%  '
$ws.Range("D10").Value = '34.61'
$ws.Range("E10").Value = '  -4.44%  '
$ws.Range("D11").Value = '0.0791'
$ws.Range("E11").Value = '  -2.35%  '
$ws.Range("E12").Value = '  +0.67%  '
$ws.Range("D13").Value = '6.75'
$ws.Range("E13").Value = '  -2.92%  '
$ws.Range("D14").Value = '2.662.35'
$ws.Range("E14").Value = '  -0.20%  '
$ws.Range("D15").Value = '15.67'
$ws.Range("E15").Value = '  +4.04%  '
$ws.Range("D16").Value = '2.316.13'
$ws.Range("E16").Value = '  +0.27%  '
$ws.Range("D17").Value = '0.801'
$ws.Range("E17").Value = '  -0.28%  '
$ws.Range("D18").Value = '42.591.58'
$ws.Range("E18").Value = '  -1.57%  '
$ws.Range("E19").Value = '  -1.79%  '
$ws.Range("D20").Value = '11.50'
$ws.Range("E20").Value = '  -3.84%  '
$ws.Range("E21").Value = '  -1.44%  '
$ws.Range("D22").Value = '68.05'
$ws.Range("E22").Value = '  +0.29%  '
$ws.Range("D23").Value = '235.08'
$ws.Range("E23").Value = '  -2.10%  '
$ws.Range("D24").Value = '1.96'
$ws.Range("E24").Value = '  -2.56%  '
$ws.Range("D25").Value = '2.51'
$ws.Range("E25").Value = '  -3.59%  '
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("D27").Value = '25.04'
$ws.Range("E27").Value = '  +0.97%  '
$ws.Range("D28").Value = '2.28'
$ws.Range("E28").Value = '  +1.69%  '
$ws.Range("D29").Value = '34.65'
$ws.Range("E29").Value = '  -4.79%  '
$ws.Range("D30").Value = '9.18'
$ws.Range("E30").Value = '  -3.91%  '
$ws.Range("D31").Value = '163.79'
$ws.Range("E31").Value = '  -0.16%  '
$ws.Range("D32").Value = '1.00'
$ws.Range("E32").Value = '  -0.04%  '
$ws.Range("D33").Value = '5.02'
$ws.Range("E33").Value = '  -4.08%  '
$ws.Range("B34").Value = 'WEMIXToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D34").Value = '2.41'
$ws.Range("E34").Value = '  -4.96%  '
$ws.Range("B35").Value = 'RenderToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D35").Value = '4.49'
$ws.Range("E35").Value = '  -0.19%  '
$ws.Range("D36").Value = '16.84'
$ws.Range("E36").Value = '  -7.71%  '
$ws.Range("D37").Value = '0.0704'
$ws.Range("E37").Value = '  -4.42%  '
$ws.Range("E38").Value = '  -4.24%  '
$ws.Range("D39").Value = '1.80'
$ws.Range("E39").Value = '  -3.61%  '
$ws.Range("D40").Value = '0.1000'
$ws.Range("E40").Value = '  -5.29%  '
$ws.Range("E41").Value = '  -3.43%  '
$ws.Range("D42").Value = '2.48'
$ws.Range("E42").Value = '  -3.22%  '
$ws.Range("D43").Value = '1.967.97'
$ws.Range("E43").Value = '  -0.79%  '
$ws.Range("D44").Value = '0.0280'
$ws.Range("E44").Value = '  -3.19%  '
$ws.Range("D45").Value = '18.49'
$ws.Range("E45").Value = '  -2.17%  '
$ws.Range("D46").Value = '10.23'
$ws.Range("E46").Value = '  +1.48%  '
$ws.Range("D47").Value = '2.88'
$ws.Range("E47").Value = '  -5.92%  '
$ws.Range("D48").Value = '55.20'
$ws.Range("E48").Value = '  -4.44%  '
$ws.Range("D49").Value = '2.528.70'
$ws.Range("E49").Value = '  -0.07%  '
$ws.Range("D50").Value = '2.83'
$ws.Range("E50").Value = '  -2.71%  '
$ws.Range("E51").Value = '  +0.23%  '
